$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "HKD" value that was in I2 (Currency column header's data row)
$ws.Range("I2").ClearContents() | Out-Null

# The active selection moved from A2 to I2
$ws.Range("I2").Select() | Out-Null

Write-Host "Applied edit: cleared I2 and updated selection to I2"
